$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sensors claimed: HiTechnicEOPD (row 20, EOPD),
# DexterPressureSensor250 (row 17, DPressure250),
# DexterPressureSensor500 (row 18, DPressure500).
# Order of entry matches the order new shared strings were created.
$ws.Range("B20").Value = "HiTechnicEOPD"
$ws.Range("B17").Value = "DexterPressureSensor250"
$ws.Range("B18").Value = "DexterPressureSensor500"
$ws.Range("F17").Value = "Pressure"
$ws.Range("F18").Value = "Pressure"
$ws.Range("F20").Value = "Distance"

# Fill remaining columns for rows 17, 18 and 20
$ws.Range("D17").Value = "Lawrie"
$ws.Range("E17").Value = "N"
$ws.Range("G17").Value = "SampleProvider"

$ws.Range("D18").Value = "Lawrie"
$ws.Range("E18").Value = "N"
$ws.Range("G18").Value = "SampleProvider"

$ws.Range("D20").Value = "Lawrie"
$ws.Range("E20").Value = "N"
$ws.Range("G20").Value = "SampleProvider"

# Update the active selection shown in the sheet view
$ws.Range("G22").Select()
